# Adds day-ahead price rows 3578-3673 (2026-01-22 / 2026-01-23 CET)
# to Sheet1, matching the source workbook's layout: column A = datetime_UTC,
# column B = datetime_CET, column C = price_eur_mwh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 3578

$aVals = @(
  46043.95833333334,
  46043.96875,
  46043.97916666666,
  46043.98958333334,
  46044,
  46044.01041666666,
  46044.02083333334,
  46044.03125,
  46044.04166666666,
  46044.05208333334,
  46044.0625,
  46044.07291666666,
  46044.08333333334,
  46044.09375,
  46044.10416666666,
  46044.11458333334,
  46044.125,
  46044.13541666666,
  46044.14583333334,
  46044.15625,
  46044.16666666666,
  46044.17708333334,
  46044.1875,
  46044.19791666666,
  46044.20833333334,
  46044.21875,
  46044.22916666666,
  46044.23958333334,
  46044.25,
  46044.26041666666,
  46044.27083333334,
  46044.28125,
  46044.29166666666,
  46044.30208333334,
  46044.3125,
  46044.32291666666,
  46044.33333333334,
  46044.34375,
  46044.35416666666,
  46044.36458333334,
  46044.375,
  46044.38541666666,
  46044.39583333334,
  46044.40625,
  46044.41666666666,
  46044.42708333334,
  46044.4375,
  46044.44791666666,
  46044.45833333334,
  46044.46875,
  46044.47916666666,
  46044.48958333334,
  46044.5,
  46044.51041666666,
  46044.52083333334,
  46044.53125,
  46044.54166666666,
  46044.55208333334,
  46044.5625,
  46044.57291666666,
  46044.58333333334,
  46044.59375,
  46044.60416666666,
  46044.61458333334,
  46044.625,
  46044.63541666666,
  46044.64583333334,
  46044.65625,
  46044.66666666666,
  46044.67708333334,
  46044.6875,
  46044.69791666666,
  46044.70833333334,
  46044.71875,
  46044.72916666666,
  46044.73958333334,
  46044.75,
  46044.76041666666,
  46044.77083333334,
  46044.78125,
  46044.79166666666,
  46044.80208333334,
  46044.8125,
  46044.82291666666,
  46044.83333333334,
  46044.84375,
  46044.85416666666,
  46044.86458333334,
  46044.875,
  46044.88541666666,
  46044.89583333334,
  46044.90625,
  46044.91666666666,
  46044.92708333334,
  46044.9375,
  46044.94791666666
)
$bVals = @(
  46044,
  46044.01041666666,
  46044.02083333334,
  46044.03125,
  46044.04166666666,
  46044.05208333334,
  46044.0625,
  46044.07291666666,
  46044.08333333334,
  46044.09375,
  46044.10416666666,
  46044.11458333334,
  46044.125,
  46044.13541666666,
  46044.14583333334,
  46044.15625,
  46044.16666666666,
  46044.17708333334,
  46044.1875,
  46044.19791666666,
  46044.20833333334,
  46044.21875,
  46044.22916666666,
  46044.23958333334,
  46044.25,
  46044.26041666666,
  46044.27083333334,
  46044.28125,
  46044.29166666666,
  46044.30208333334,
  46044.3125,
  46044.32291666666,
  46044.33333333334,
  46044.34375,
  46044.35416666666,
  46044.36458333334,
  46044.375,
  46044.38541666666,
  46044.39583333334,
  46044.40625,
  46044.41666666666,
  46044.42708333334,
  46044.4375,
  46044.44791666666,
  46044.45833333334,
  46044.46875,
  46044.47916666666,
  46044.48958333334,
  46044.5,
  46044.51041666666,
  46044.52083333334,
  46044.53125,
  46044.54166666666,
  46044.55208333334,
  46044.5625,
  46044.57291666666,
  46044.58333333334,
  46044.59375,
  46044.60416666666,
  46044.61458333334,
  46044.625,
  46044.63541666666,
  46044.64583333334,
  46044.65625,
  46044.66666666666,
  46044.67708333334,
  46044.6875,
  46044.69791666666,
  46044.70833333334,
  46044.71875,
  46044.72916666666,
  46044.73958333334,
  46044.75,
  46044.76041666666,
  46044.77083333334,
  46044.78125,
  46044.79166666666,
  46044.80208333334,
  46044.8125,
  46044.82291666666,
  46044.83333333334,
  46044.84375,
  46044.85416666666,
  46044.86458333334,
  46044.875,
  46044.88541666666,
  46044.89583333334,
  46044.90625,
  46044.91666666666,
  46044.92708333334,
  46044.9375,
  46044.94791666666,
  46044.95833333334,
  46044.96875,
  46044.97916666666,
  46044.98958333334
)
$cVals = @(
  88.03,
  88.42,
  85.73,
  80.92,
  89.89,
  88.53,
  87.72,
  86,
  90.18000000000001,
  90.09,
  89.98,
  84.56999999999999,
  85.43000000000001,
  83.65000000000001,
  83.47,
  83.81999999999999,
  83.86,
  83.61,
  85.06,
  85.95999999999999,
  84.86,
  82.98,
  84.89,
  91.8,
  85.23,
  88.73,
  101.37,
  111.26,
  97.34999999999999,
  104.82,
  106.4,
  119.13,
  135.01,
  143.18,
  128.97,
  122.18,
  127.5,
  115.67,
  112.61,
  100,
  109.02,
  98.59,
  98.33,
  94.34,
  91.5,
  89.73999999999999,
  89.94,
  87.45,
  90.73,
  87.01000000000001,
  86.53,
  85.88,
  89.94,
  90.76000000000001,
  90.06,
  90.17,
  90.67,
  97.09,
  96.87,
  101.49,
  95.83,
  101.88,
  108.25,
  113.86,
  110.07,
  123.16,
  136.94,
  127.61,
  123.63,
  131.92,
  136.64,
  134.98,
  126.39,
  127.39,
  127.42,
  130,
  128.65,
  120,
  122.67,
  112.81,
  115.91,
  108.31,
  108.1,
  106.14,
  120.33,
  108.04,
  101.84,
  97.44,
  102.68,
  101.91,
  108.95,
  98.63,
  102.67,
  96.56,
  99.45,
  90.45999999999999
)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $aVals[$i]
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
}

$lastRow = $startRow + $aVals.Length - 1
$ws.Range("A" + $startRow + ":B" + $lastRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"
